# IPC_Evolucion_conApertura.xlsx - "nuevo cambio en IPC calculado"
# Adds "Indice_Final" (echo of index columns C/D/E into G/H/I) and
# "Var Anual" (year-over-year variation, moved from G/H/I into J/K/L)
# plus a brand-new row (row 15) with the latest period's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------
# Row 1 (top header band): G1:I1 -> "Indice_Final", J1:L1 -> "Var Anual"
# ---------------------------------------------------------------
$ws.Range("G1").Value = "Indice_Final"
$ws.Range("H1").Value = "Indice_Final"
$ws.Range("I1").Value = "Indice_Final"

$ws.Range("J1").Value = "Var Anual"
$ws.Range("J1").HorizontalAlignment = -4108
$ws.Range("K1").Value = "Var Anual"
$ws.Range("K1").HorizontalAlignment = -4108
$ws.Range("L1").Value = "Var Anual"
$ws.Range("L1").HorizontalAlignment = -4108

# ---------------------------------------------------------------
# Row 2 (sub header): copy the G2:I2 formatting (fill) into J2:L2
# and give them the same labels 00 / 01 / 02
# ---------------------------------------------------------------
$ws.Range("G2").Copy() | Out-Null
$ws.Range("J2").PasteSpecial(-4122) | Out-Null
$ws.Range("J2").Value = "00"

$ws.Range("H2").Copy() | Out-Null
$ws.Range("K2").PasteSpecial(-4122) | Out-Null
$ws.Range("K2").Value = "01"

$ws.Range("I2").Copy() | Out-Null
$ws.Range("L2").PasteSpecial(-4122) | Out-Null
$ws.Range("L2").Value = "02"

$ws.Application.CutCopyMode = 0

# ---------------------------------------------------------------
# Row 3: add F3:I3 echoes of A3/C3/D3/E3 (brand new cells)
# ---------------------------------------------------------------
$ws.Range("F3").Formula = "=A3"
$ws.Range("F3").NumberFormat = "mmm-yy"
$ws.Range("G3").Formula = "=C3"
$ws.Range("H3").Formula = "=D3"
$ws.Range("I3").Formula = "=E3"
$ws.Range("G3:I3").NumberFormat = "0.00"

# ---------------------------------------------------------------
# Rows 4-11: move the existing "year over year" formulas that live in
# G:I today over to J:K:L, then replace G:I with simple echoes of C:D:E
# ---------------------------------------------------------------
for ($r = 4; $r -le 11; $r++) {
    $gFormula = $ws.Range("G$r").Formula
    $hFormula = $ws.Range("H$r").Formula
    $iFormula = $ws.Range("I$r").Formula

    $ws.Range("J$r").Formula = $gFormula
    $ws.Range("K$r").Formula = $hFormula
    $ws.Range("L$r").Formula = $iFormula
    $ws.Range("J$r`:L$r").NumberFormat = "0.00%"

    $ws.Range("G$r").Formula = "=C$r"
    $ws.Range("H$r").Formula = "=D$r"
    $ws.Range("I$r").Formula = "=E$r"
    $ws.Range("G$r`:I$r").NumberFormat = "0.00"
}

# ---------------------------------------------------------------
# Row 12: brand new F12:I12 cells (echoes only, no J:K:L here)
# ---------------------------------------------------------------
$ws.Range("F12").Formula = "=A12"
$ws.Range("F12").NumberFormat = "mmm-yy"
$ws.Range("G12").Formula = "=C12"
$ws.Range("H12").Formula = "=D12"
$ws.Range("I12").Formula = "=E12"
$ws.Range("G12:I12").NumberFormat = "0.00"

# ---------------------------------------------------------------
# Row 13: blank G13:I13 cells, just formatted as 0.00
# ---------------------------------------------------------------
$ws.Range("G13").NumberFormat = "0.00"
$ws.Range("H13").NumberFormat = "0.00"
$ws.Range("I13").NumberFormat = "0.00"

# ---------------------------------------------------------------
# Row 14: small data correction + new formulas in G:I, moved J:K:L formulas
# ---------------------------------------------------------------
$ws.Range("C14").Value = 99.465010000000007

$ws.Range("J14").Formula = "=(1+(C12/C10-1))*(1+(C14/C13-1))-1"
$ws.Range("K14").Formula = "=(1+(D12/D10-1))*(1+(D14/D13-1))-1"
$ws.Range("L14").Formula = "=(1+(E12/E10-1))*(1+(E14/E13-1))-1"
$ws.Range("J14:L14").NumberFormat = "0.00%"

$ws.Range("G14").Formula = "=C14/C`$13*C`$12"
$ws.Range("H14").Formula = "=D14/D`$13*D`$12"
$ws.Range("I14").Formula = "=E14/E`$13*E`$12"
$ws.Range("G14:I14").NumberFormat = "0.00"

# ---------------------------------------------------------------
# Row 15 (new): full new row of data
# ---------------------------------------------------------------
$ws.Range("A15").Value = 45107
$ws.Range("A15").NumberFormat = "mmm-yy"
$ws.Range("B15").Value = 44865
$ws.Range("B15").NumberFormat = "mmm-yy"
$ws.Range("C15").Value = 103.22
$ws.Range("C15").NumberFormat = "0.00"
$ws.Range("D15").Value = 106.10509999999999
$ws.Range("D15").NumberFormat = "0.00"
$ws.Range("E15").Value = 104.84
$ws.Range("E15").NumberFormat = "0.00"

$ws.Range("F15").Formula = "=A15"
$ws.Range("F15").NumberFormat = "mmm-yy"

$ws.Range("G15").Formula = "=C15/C`$13*C`$12"
$ws.Range("H15").Formula = "=D15/D`$13*D`$12"
$ws.Range("I15").Formula = "=E15/E`$13*E`$12"
$ws.Range("G15:I15").NumberFormat = "0.00"

$ws.Range("J15").Formula = "=G15/G11-1"
$ws.Range("K15").Formula = "=H15/H11-1"
$ws.Range("L15").Formula = "=I15/I11-1"
$ws.Range("J15:L15").NumberFormat = "0.00%"

$ws.Range("M15").NumberFormat = "0.00%"
$ws.Range("N15").NumberFormat = "0.00%"
$ws.Range("O15").NumberFormat = "0.00%"

# ---------------------------------------------------------------
# Selection / view bookkeeping to match the target workbook state
# ---------------------------------------------------------------
$ws.Range("H22").Select() | Out-Null

$wb.Save()
